# SysMetaMutability.xlsx — update ControlledBy / Method columns so that the
# MN (Member Node) controls each system-metadata field except the replica
# entries, which remain CN-controlled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# replica: Method gains explicit call-parentheses (still CN-controlled)
$ws.Range("E18").Value = "CNReplication.updateReplicationMetadata()"

# authoritativeMemberNode: ControlledBy CN -> MN/CN, ModifiableBy
$ws.Range("C17").Value = "MN/CN"
$ws.Range("F17").Value = "Someone with access to MNs/CNs"

# archived: Method
$ws.Range("E13").Value = "MNCore.archive()"

# obsoletes: Method
$ws.Range("E11").Value = "MNStorage.update()"

# obsoletedBy: Method
$ws.Range("E12").Value = "MNCore.setObsoletedBy(), MNStorage.update()"

# accessPolicy: Method (also taller row to fit the longer text)
$ws.Range("E9").Value = "manual (Tier 1), MNAuthorization.setAccessPolicy(), MNStorage.update ()(all must call CNAuthorization.systemMetadataChanged())"
$ws.Rows.Item(9).RowHeight = 45

# rightsHolder: Method
$ws.Range("E8").Value = "MNAuthorization.setRightsHolder()"

# serialVersion: ControlledBy CN -> MN, ModifiableBy CN-service-subject -> MN-service-subject
$ws.Range("F2").Value = "MN-service-subject"
$ws.Range("C2").Value = "MN"

# dateSysMetadataModified: ControlledBy CN -> MN
$ws.Range("C15").Value = "MN"

# originMemberNode: ControlledBy CN -> MN
$ws.Range("C16").Value = "MN"

# Restore the selection to match the saved view (F4)
$ws.Range("F4").Select()
